$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.144984483718872
$ws.Range("B1").Value = 2.353566884994507
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 1.801346778869629
$ws.Range("E1").Value = 1.208109259605408
